# Append the latest EUR->ARS quote as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

# Force column A's cell to text first so the date-looking string
# ("2025-09-22") is stored as a literal string instead of being
# auto-converted into a serial date number, then restore the default
# "Normal" style so no extra formatting is left on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-22"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "21:21:38"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,732.1025"
